# Updated cryptos list on Fri Oct  6 11:59:33 UTC 2023 with GitHub Actions
#
# Helper: write a value into a cell while keeping it a plain text string
# (these "Price" figures use '.' as both a decimal AND thousands separator,
# e.g. "212.07" or "1.634.07", and Excel would otherwise auto-coerce the
# single-dot ones into numbers on assignment). Forcing the NumberFormat to
# Text before the write, then resetting the style back to Normal afterwards,
# keeps the stored value a string without leaving a lasting number format.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.728.92"
$ws.Range("E2").Value = "  -0.08%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.634.07"
$ws.Range("E3").Value = "  -0.24%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "212.07"
$ws.Range("E5").Value = "  -0.22%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.31%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - Solana
Set-TextValue $ws.Range("D8") "23.21"
$ws.Range("E8").Value = "  -0.15%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.91%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.22%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -3.06%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.865.06"
$ws.Range("E12").Value = "  -0.29%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.638.40"
$ws.Range("E13").Value = "  -0.17%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -0.30%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -1.41%  "

# Row 16 - Litecoin
Set-TextValue $ws.Range("D16") "65.17"

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "27.676.56"
$ws.Range("E17").Value = "  -0.21%  "

# Row 18 - BitcoinCash
Set-TextValue $ws.Range("D18") "229.40"
$ws.Range("E18").Value = "  -0.20%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  -0.38%  "

# Row 20 - Chainlink
Set-TextValue $ws.Range("D20") "7.57"
$ws.Range("E20").Value = "  -1.68%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.05%  "

# Row 22 - Avalanche
Set-TextValue $ws.Range("D22") "10.66"
$ws.Range("E22").Value = "  +4.51%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +0.73%  "

# Row 25 - Monero
Set-TextValue $ws.Range("D25") "149.08"
$ws.Range("E25").Value = "  -1.39%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  -1.05%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -1.00%  "

# Row 28 - EthereumClassic
Set-TextValue $ws.Range("D28") "15.58"
$ws.Range("E28").Value = "  -0.12%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.01%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.28%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -0.95%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -0.72%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.471.48"
$ws.Range("E33").Value = "  +0.01%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  -1.14%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -1.49%  "

# Row 37 - TrustWalletToken
Set-TextValue $ws.Range("D37") "0.933"
$ws.Range("E37").Value = "  +2.04%  "

# Rows 38/39 - ImmutableX and ARBITRUM swap places (ARBITRUM now ranks 38th)
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D38") "0.877"
$ws.Range("E38").Value = "  -0.46%  "

$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D39") "0.559"
$ws.Range("E39").Value = "  -1.60%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +0.01%  "

# Row 41 - WEMIXToken
Set-TextValue $ws.Range("D41") "1.02"
$ws.Range("E41").Value = "  +0.51%  "

# Row 42 - Aave
$ws.Range("E42").Value = "  -1.54%  "

# Row 43 - mCoin
$ws.Range("E43").Value = "  +0.90%  "

# Row 44 - MXToken
$ws.Range("E44").Value = "  -0.95%  "

# Row 45 - FraxShare
Set-TextValue $ws.Range("D45") "5.37"
$ws.Range("E45").Value = "  -4.08%  "

# Row 46 - RocketPoolETH
$ws.Range("D46").Value = "1.774.74"
$ws.Range("E46").Value = "  -0.35%  "

# Row 47 - RenderToken
Set-TextValue $ws.Range("D47") "1.74"
$ws.Range("E47").Value = "  +1.76%  "

# Row 48 - Quant
Set-TextValue $ws.Range("D48") "87.63"
$ws.Range("E48").Value = "  +0.80%  "

# Row 49 - BabyDogeCoin
$ws.Range("E49").Value = "  -0.68%  "

# Row 50 - Algorand
Set-TextValue $ws.Range("D50") "0.0992"
$ws.Range("E50").Value = "  -0.12%  "

# Row 51 - EnergySwap
Set-TextValue $ws.Range("D51") "7.69"
$ws.Range("E51").Value = "  -1.81%  "
